$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KDP")

# Row 4: Inventory
$ws.Range("B4").Value = 841000000.0
$ws.Range("C4").Value = 817000000.0
$ws.Range("D4").Value = 876000000.0
$ws.Range("E4").Value = 797000000.0
$ws.Range("F4").Value = 732000000.0

# Row 15: Accounts Payable
$ws.Range("B15").Value = 3871000000.0
$ws.Range("C15").Value = 3740000000.0
$ws.Range("D15").Value = 3517000000.0
$ws.Range("E15").Value = 3377000000.0
$ws.Range("F15").Value = 3238000000.0

# Row 25: Long Term Tax Liability (Deferred)
$ws.Range("B25").Value = 5981000000.0
$ws.Range("C25").Value = 5948000000.0
$ws.Range("D25").Value = 5914000000.0
$ws.Range("E25").Value = 5893000000.0
$ws.Range("F25").Value = 5888000000.0

# Row 39: Net Debt
$ws.Range("G39").Value = 14750000000.0

# Row 40: Total Debt
$ws.Range("G40").Value = 14851000000.0
